$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = "5 Spinach Linguine"
$ws.Range("B3").Value = "5 SPLING 4E"
$ws.Range("D8").Select()
